$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 rows right before the "Azure" row (currently row 62).
# This makes room for two new single-row skill entries ("dbt" and "Tableau"),
# each occupying the same 3-row block pattern used by the other skills
# (Azure, AWS, Snowflake, Flink, ...), and pushes Azure/AWS down by 6 rows.
$ws.Rows("62:67").Insert()

# Fill in the two new skill rows.
$ws.Range("C62").Value = "dbt"
$ws.Range("C65").Value = "Tableau"

# Insert 8 more rows further down, inside the blank gap that exists between
# AWS (now at row 71) and DBMS (currently at row 82), pushing DBMS/OS/CN
# down by a further 8 rows (82->96, 85->99, 88->102), matching the target
# layout, without disturbing Azure/AWS above.
$ws.Rows("75:82").Insert()

# Update the active selection to match the saved view state.
$ws.Range("D31").Select()
